# "Generate Report for Handoff"
# Updates the status/handoff info for file b.md (row 3) across the
# Overview, zh-cn and de-de sheets: it is now ready for handoff again
# (handback was stale), with a new handoff xliff + timestamp and an
# error detail message.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 corresponds to b.md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-05 02:45:22"

# --- zh-cn sheet: row 3 corresponds to b.md ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"

# "False" would be auto-detected as a Boolean by a plain .Value assignment,
# so force it through a text formula and bake the formula back down to a
# literal value via copy / paste-special.
$wsZhCn.Range("F3").Formula = "=""False"""
$wsZhCn.Range("F3").Copy()
$wsZhCn.Range("F3").PasteSpecial(-4163)

$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-05 02:45:11"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/c556bfcb02303818985b426766897ec6b7a5faf3/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/e6b46ebf2a073f94f42b6f02f155efaa17aa1248/e2e/b.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet: row 3 corresponds to b.md ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"

$wsDeDe.Range("F3").Formula = "=""False"""
$wsDeDe.Range("F3").Copy()
$wsDeDe.Range("F3").PasteSpecial(-4163)

$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-05 02:45:22"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/c556bfcb02303818985b426766897ec6b7a5faf3/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/e6b46ebf2a073f94f42b6f02f155efaa17aa1248/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
